$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Result, angry, happy, neutral, sad)
$data = @{
    2  = @("correct",   1,                  0,                  0,                  0)
    3  = @("correct",   1,                  0,                  0,                  0)
    4  = @("correct",   1,                  0,                  0,                  0)
    5  = @("correct",   1,                  0,                  0,                  0)
    6  = @("correct",   1,                  0,                  0,                  0)
    7  = @("correct",   1,                  0,                  0,                  0)
    8  = @("incorrect", 0.3386495073309164, 0.3306152359259172, 0,                  0.3307352567431665)
    9  = @("incorrect", 1,                  0,                  0,                  0)
    10 = @("correct",   0,                  1,                  0,                  0)
    11 = @("incorrect", 1,                  0,                  0,                  0)
    12 = @("incorrect", 0.6653830943429145, 0,                  0,                  0.3346169056570855)
    13 = @("correct",   0,                  1,                  0,                  0)
    15 = @("correct",   0,                  0,                  1,                  0)
    16 = @("correct",   0,                  0,                  1,                  0)
    17 = @("correct",   0,                  0,                  1,                  0)
    18 = @("correct",   0,                  0,                  1,                  0)
    19 = @("correct",   0,                  0,                  1,                  0)
    20 = @("incorrect", 1,                  0,                  0,                  0)
    21 = @("correct",   0,                  0,                  0,                  1)
    22 = @("incorrect", 0.6685935153207773, 0,                  0.3314064846792227, 0)
    23 = @("correct",   0,                  0,                  0,                  1)
    24 = @("incorrect", 0.6671008803626376, 0.3328991196373624, 0,                  0)
    25 = @("incorrect", 0.671362379512858,  0,                  0.3286376204871418, 0)
    26 = @("incorrect", 0,                  0,                  1,                  0)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}
